$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column E (Obrigatorio) for rows 2 through 9 from "N" to "S"
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 5).Value = "S"
}
